$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: remove "CM - Cost Structure" / "test cost structure" test case, keep formatting
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()

# Row 4: remove "CM - Pricing" / "pricing" test case, keep formatting
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# New blank styled row 5 cell (copy formatting from an existing styled helper cell)
$ws.Range("G4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = ""
$excel.CutCopyMode = $false

# Row 6: remove old "Top 100 Aging sku" test case row entirely
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Clear()

# Update selection / scroll position shown in the saved view
$ws.Range("E11").Select()
